# 03项目计划表.xlsx - add new "2018.10.15" planning section (rows 41-48)
# and a fresh summary block (rows 49-50), while turning the former
# placeholder summary block (rows 39-40) into the real write-up for the
# previous ("2018.10.11") section.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Fill in the summary text that used to be an empty "总结：" placeholder
#    (rows 39:40, merged) with the actual write-up.
# ---------------------------------------------------------------------
$ws.Range("A39").Value = "总结：经过小组讨论之后觉得之前的界面确实是有点简陋，美感不足，无法吸引用户，觉得一个好的界面不仅可以吸引用户，还可以让我们的价值得到提升，苦寻之下我们终于有了创作的源泉，于是我们墨刀设计我们小组专属的app界面，我们坚信天道酬勤，苦心人天不负，百二秦关终属楚"

# ---------------------------------------------------------------------
# 2. Insert 10 fresh rows (41-50) below the existing block, and build
#    the new "日期：2018.10.15第七周周一" table there, followed by a
#    brand-new (still empty) "总结：" placeholder.
# ---------------------------------------------------------------------
$ws.Rows.Item(41).Resize(10).Insert()

# -- Title row (41), copy look & feel of the previous title row (31) --
$ws.Range("A31:D31").Copy()
$ws.Range("A41:D41").PasteSpecial(-4122) | Out-Null
$ws.Range("A41").Value = "日期：2018.10.15第七周周一"
$ws.Range("A41:D41").Merge() | Out-Null

# -- Header row (42), copy from the previous header row (32) --
$ws.Range("A32:D32").Copy()
$ws.Range("A42:D42").PasteSpecial(-4122) | Out-Null
$ws.Range("A42").Value = "组员"
$ws.Range("B42").Value = "计划内容"
$ws.Range("C42").Value = "完成情况"
$ws.Range("D42").Value = "备注"

# -- Data rows (43-48), copy from the previous data rows (33-38) --
$ws.Range("A33:D38").Copy()
$ws.Range("A43:D48").PasteSpecial(-4122) | Out-Null

$ws.Range("A43").Value = "陈柯赞"
$ws.Range("B43").Value = "数据库用户表chat_user设计及E-R图"
$ws.Range("C43").Value = ""

$ws.Range("A44").Value = "黎安生"
$ws.Range("B44").Value = "数据库管理员表chat_administrator设计及E-R图"
$ws.Range("C44").Value = ""

$ws.Range("A45").Value = "王智永"
$ws.Range("B45").Value = "数据库群组表chat_group设计及E-R图"
$ws.Range("C45").Value = ""

$ws.Range("A46").Value = "郑海文"
$ws.Range("B46").Value = "数据库chat_group_user设计及E-R图"
$ws.Range("C46").Value = ""

$ws.Range("A47").Value = "赵华亮"
$ws.Range("B47").Value = "数据库聊天表chat_message设计及E-R图"
$ws.Range("C47").Value = ""

$ws.Range("A48").Value = "叶田"
$ws.Range("B48").Value = "数据库定位表chat_user_gps"
$ws.Range("C48").Value = ""

# -- Remarks column: single note merged across D43:D48 --
$ws.Range("D43").Value = "我们使用startUml画数据库E-R图"
$ws.Range("D44").Value = ""
$ws.Range("D45").Value = ""
$ws.Range("D46").Value = ""
$ws.Range("D47").Value = ""
$ws.Range("D48").Value = ""
$ws.Range("D43:D48").Merge() | Out-Null
$ws.Range("D43:D48").HorizontalAlignment = -4108 # xlCenter
$ws.Range("D43:D48").VerticalAlignment = -4108   # xlCenter
$ws.Range("D43:D48").Borders.Item(7).LineStyle = 1  # xlEdgeLeft
$ws.Range("D43:D48").Borders.Item(10).LineStyle = 1 # xlEdgeRight
$ws.Range("D43").Borders.Item(8).LineStyle = 1      # xlEdgeTop
$ws.Range("D48").Borders.Item(9).LineStyle = 1      # xlEdgeBottom

# -- New empty summary placeholder (49-50), copy from rows 39-40 --
$ws.Range("A39:D40").Copy()
$ws.Range("A49:D50").PasteSpecial(-4122) | Out-Null
$ws.Range("A49").Value = "总结："
$ws.Range("B49").Value = ""
$ws.Range("C49").Value = ""
$ws.Range("D49").Value = ""
$ws.Range("A50").Value = ""
$ws.Range("B50").Value = ""
$ws.Range("C50").Value = ""
$ws.Range("D50").Value = ""
$ws.Range("A49:D50").Merge() | Out-Null

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3. Column D needs to be much wider now that it holds a paragraph of
#    remarks, and the view should scroll/select the new block.
# ---------------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 89.5

$ws.Application.ActiveWindow.ScrollRow = 19
$ws.Range("D43:D48").Select() | Out-Null
